$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Extend header row (row 1) formatting to the new columns H:N by
#     copying the style from an existing header cell (G1 -> style s="1") ---
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# --- Extend data row (row 2) formatting to the new columns H:N by
#     copying the style from an existing data cell (G2 -> style s="2") ---
$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row 1 (header) values ---
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 (data) values ---
$ws.Range("A2").Value = 110
$ws.Range("B2").Value = "擔保借款"
$ws.Range("C2").Value = "張嘉郡"
$ws.Range("D2").Value = "台新銀行苓雅分行高雄市前鎮區中山二路"
$ws.Range("E2").Value = 1221524
$ws.Range("F2").Value = "100年04月08日"
$ws.Range("G2").Value = "購車"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"

# J2 needs the literal text "2011-11-18" - assigning that string directly
# would get auto-recognised as a date (like real Excel date entry), turning
# it into a date serial instead of the source text. Format the cell as Text
# first so the value is stored verbatim as a string.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-11-18"

$ws.Range("K2").Value = "張嘉郡"
$ws.Range("L2").Value = 1719
$ws.Range("M2").Value = "tmp77961"
$ws.Range("N2").Value = 110
